# Arbeitszeit.xlsx - update hours for B21 and move selection to B22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tracked hours for the last entry from 8 to 10
$ws.Range("B21").Value = 10

# Force recalculation so the SUM formula in E2 reflects the new value
$excel.Calculate()

# Move the active selection to B22 (next empty row), matching the diff
$ws.Range("B22").Select()
